# Update countries & provincias Spain
# Refresh the "Pais" sheet with the newer COVID-19 snapshot:
#  - update the "last updated" timestamp
#  - refresh per-country metrics (Casos totales, Nuevos casos, Casos activos,
#    Recuperados, Casos criticos, Muertes hoy, Muertes)
#  - Colombia overtakes Mexico, Argentina overtakes Iran, Montenegro overtakes
#    Nicaragua and Polinesia Francesa overtakes Islas Feroe in the ranking,
#    so those row pairs swap country/data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 27 de Agosto de 2020 a las 01:47'

$ws.Range("B4").Value = 5997927
$ws.Range("C4").Value = 42199
$ws.Range("D4").Value = 3295346
$ws.Range("E4").Value = 2518988
$ws.Range("G4").Value = 1229
$ws.Range("H4").Value = 183593

$ws.Range("B5").Value = 3722004
$ws.Range("C5").Value = 47828
$ws.Range("E5").Value = 695400
$ws.Range("G5").Value = 1090
$ws.Range("H5").Value = 117756

$ws.Range("A10").Value = 'Colombia'
$ws.Range("B10").Value = 572270
$ws.Range("C10").Value = 10142
$ws.Range("D10").Value = 407121
$ws.Range("E10").Value = 146965
$ws.Range("G10").Value = 295
$ws.Range("H10").Value = 18184

$ws.Range("A11").Value = 'Mexico'
$ws.Range("B11").Value = 568621
$ws.Range("C11").Value = 4916
$ws.Range("D11").Value = 393101
$ws.Range("E11").Value = 114070
$ws.Range("G11").Value = 650
$ws.Range("H11").Value = 61450

$ws.Range("A14").Value = 'Argentina'
$ws.Range("B14").Value = 370188
$ws.Range("C14").Value = 10550
$ws.Range("D14").Value = 268801
$ws.Range("E14").Value = 93548
$ws.Range("G14").Value = 276
$ws.Range("H14").Value = 7839

$ws.Range("A15").Value = 'Iran'
$ws.Range("B15").Value = 365606
$ws.Range("C15").Value = 2243
$ws.Range("D15").Value = 314870
$ws.Range("E15").Value = 29716
$ws.Range("G15").Value = 119
$ws.Range("H15").Value = 21020

$ws.Range("D23").Value = 214233
$ws.Range("E23").Value = 15415

$ws.Range("B36").Value = 89082
$ws.Range("C36").Value = 701
$ws.Range("D36").Value = 63595
$ws.Range("E36").Value = 23555
$ws.Range("G36").Value = 13
$ws.Range("H36").Value = 1932

$ws.Range("B53").Value = 53021
$ws.Range("C53").Value = 221
$ws.Range("D53").Value = 40281
$ws.Range("E53").Value = 11730
$ws.Range("G53").Value = 3
$ws.Range("H53").Value = 1010

$ws.Range("B74").Value = 22951
$ws.Range("C74").Value = 403
$ws.Range("D74").Value = 16954
$ws.Range("E74").Value = 5579

$ws.Range("A113").Value = 'Montenegro'
$ws.Range("B113").Value = 4499
$ws.Range("C113").Value = 55
$ws.Range("D113").Value = 3558
$ws.Range("E113").Value = 853
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = 88

$ws.Range("A114").Value = 'Nicaragua'
$ws.Range("B114").Value = 4494
$ws.Range("D114").Value = 2913
$ws.Range("E114").Value = 1444
$ws.Range("H114").Value = 137

$ws.Range("B118").Value = 3724
$ws.Range("C118").Value = 26
$ws.Range("D118").Value = 2863
$ws.Range("E118").Value = 799
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 62

$ws.Range("B126").Value = 3206
$ws.Range("C126").Value = 137
$ws.Range("D126").Value = 1476
$ws.Range("E126").Value = 1659

$ws.Range("B149").Value = 1543
$ws.Range("C149").Value = 7
$ws.Range("D149").Value = 1322
$ws.Range("E149").Value = 178

$ws.Range("B152").Value = 1411
$ws.Range("C152").Value = 159
$ws.Range("E152").Value = 1204

$ws.Range("A178").Value = 'Polinesia Francesa'
$ws.Range("B178").Value = 415
$ws.Range("C178").Value = 43
$ws.Range("D178").Value = 202
$ws.Range("E178").Value = 213

$ws.Range("A179").Value = 'Islas Feroe'
$ws.Range("B179").Value = 411
$ws.Range("D179").Value = 357
$ws.Range("E179").Value = 54
